$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2's cells (value + formatting) into the new row 4, and move
# the password value from B2 down into B3 - using Copy-to-destination so both
# the value and the existing cell style/formatting travel together.
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("B2").Copy($ws.Range("B4"))
$ws.Range("B2").Copy($ws.Range("B3"))

# Clear out B2's value now that it has been relocated to B3, keeping its
# existing formatting intact.
$ws.Range("B2").ClearContents()

# Update the active selection to B2, as in the target workbook.
$ws.Range("B2").Select()
